# Update the worksheet date and the five rows of division problems.
$d = $word.ActiveDocument

# 1. Update the date/weekday heading paragraph.
$d.Content.Find.Execute("2025-04-02 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-04-03 Thursday", 2)

# 2. Update each division-problem cell in the first table, row by row.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "879÷6="
$t.Cell(1, 2).Range.Text = "726÷9="
$t.Cell(1, 3).Range.Text = "770÷5="
$t.Cell(1, 4).Range.Text = "417÷8="
$t.Cell(1, 5).Range.Text = "414÷5="

$t.Cell(5, 1).Range.Text = "640÷5="
$t.Cell(5, 2).Range.Text = "580÷2="
$t.Cell(5, 3).Range.Text = "105÷3="
$t.Cell(5, 4).Range.Text = "520÷4="
$t.Cell(5, 5).Range.Text = "834÷9="

$t.Cell(9, 1).Range.Text = "956÷4="
$t.Cell(9, 2).Range.Text = "558÷3="
$t.Cell(9, 3).Range.Text = "542÷3="
$t.Cell(9, 4).Range.Text = "836÷3="
$t.Cell(9, 5).Range.Text = "763÷4="

$t.Cell(13, 1).Range.Text = "522÷9="
$t.Cell(13, 2).Range.Text = "763÷9="
$t.Cell(13, 3).Range.Text = "166÷9="
$t.Cell(13, 4).Range.Text = "696÷3="
$t.Cell(13, 5).Range.Text = "614÷2="

$t.Cell(17, 1).Range.Text = "787÷8="
$t.Cell(17, 2).Range.Text = "434÷6="
$t.Cell(17, 3).Range.Text = "801÷3="
$t.Cell(17, 4).Range.Text = "740÷6="
$t.Cell(17, 5).Range.Text = "412÷6="
